$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to hold the literal string value (never auto-detected
    # as a number/date), then restore the default "Normal" style so no
    # stray number-format style is left behind on cells that were plain.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '35.029.25'
Set-TextValue $ws.Range('E2') '  +0.97%  '

Set-TextValue $ws.Range('D3') '1.852.09'
Set-TextValue $ws.Range('E3') '  +2.28%  '

Set-TextValue $ws.Range('E4') '  +0.19%  '

Set-TextValue $ws.Range('D5') '236.76'
Set-TextValue $ws.Range('E5') '  +3.01%  '

Set-TextValue $ws.Range('E6') '  +0.89%  '

Set-TextValue $ws.Range('E7') '  +0.12%  '

Set-TextValue $ws.Range('D8') '42.36'
Set-TextValue $ws.Range('E8') '  +6.78%  '

Set-TextValue $ws.Range('E9') '  +2.37%  '

Set-TextValue $ws.Range('E10') '  +2.00%  '

Set-TextValue $ws.Range('D11') '0.0993'
Set-TextValue $ws.Range('E11') '  +0.61%  '

Set-TextValue $ws.Range('D12') '2.120.61'
Set-TextValue $ws.Range('E12') '  +2.33%  '

Set-TextValue $ws.Range('B13') 'WrappedEther'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D13') '1.861.15'
Set-TextValue $ws.Range('E13') '  -3.74%  '

Set-TextValue $ws.Range('B14') 'Chainlink'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D14') '11.42'
Set-TextValue $ws.Range('E14') '  +2.08%  '

Set-TextValue $ws.Range('B15') 'Polkadot'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D15') '4.80'
Set-TextValue $ws.Range('E15') '  +4.99%  '

Set-TextValue $ws.Range('D16') '0.676'
Set-TextValue $ws.Range('E16') '  +1.78%  '

Set-TextValue $ws.Range('D17') '35.014.52'
Set-TextValue $ws.Range('E17') '  +1.30%  '

Set-TextValue $ws.Range('E18') '  +1.44%  '

Set-TextValue $ws.Range('D19') '0.0₃0795'
Set-TextValue $ws.Range('E19') '  +1.89%  '

Set-TextValue $ws.Range('D20') '240.57'
Set-TextValue $ws.Range('E20') '  +0.45%  '

Set-TextValue $ws.Range('D21') '12.18'
Set-TextValue $ws.Range('E21') '  +2.99%  '

Set-TextValue $ws.Range('E22') '  +3.17%  '

Set-TextValue $ws.Range('E23') '  -0.04%  '

Set-TextValue $ws.Range('E24') '  +1.41%  '

Set-TextValue $ws.Range('D25') '170.54'
Set-TextValue $ws.Range('E25') '  -1.53%  '

Set-TextValue $ws.Range('D26') '1.87'
Set-TextValue $ws.Range('E26') '  +24.49%  '

Set-TextValue $ws.Range('D27') '7.98'
Set-TextValue $ws.Range('E27') '  +3.46%  '

Set-TextValue $ws.Range('D28') '17.66'
Set-TextValue $ws.Range('E28') '  +2.27%  '

Set-TextValue $ws.Range('E29') '  +0.20%  '

Set-TextValue $ws.Range('B30') 'BinanceUSD'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D30') '1.01'
Set-TextValue $ws.Range('E30') '  +0.24%  '

Set-TextValue $ws.Range('B31') 'Hedera'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D31') '0.0558'
Set-TextValue $ws.Range('E31') '  +2.49%  '

Set-TextValue $ws.Range('D32') '4.01'
Set-TextValue $ws.Range('E32') '  +0.76%  '

Set-TextValue $ws.Range('E33') '  +3.20%  '

Set-TextValue $ws.Range('E34') '  +23.90%  '

Set-TextValue $ws.Range('E35') '  +11.77%  '

Set-TextValue $ws.Range('E36') '  +7.86%  '

Set-TextValue $ws.Range('D37') '0.780'
Set-TextValue $ws.Range('E37') '  +13.78%  '

Set-TextValue $ws.Range('E38') '  +11.27%  '

Set-TextValue $ws.Range('E39') '  +5.96%  '

Set-TextValue $ws.Range('D40') '90.64'
Set-TextValue $ws.Range('E40') '  -0.20%  '

Set-TextValue $ws.Range('D41') '1.350.20'
Set-TextValue $ws.Range('E41') '  +1.42%  '

Set-TextValue $ws.Range('E42') '  +3.86%  '

Set-TextValue $ws.Range('D43') '2.33'
Set-TextValue $ws.Range('E43') '  +3.42%  '

Set-TextValue $ws.Range('D44') '12.72'
Set-TextValue $ws.Range('E44') '  +51.54%  '

Set-TextValue $ws.Range('E45') '  -0.13%  '

Set-TextValue $ws.Range('E46') '  +6.53%  '

Set-TextValue $ws.Range('E47') '  +0.01%  '

Set-TextValue $ws.Range('E48') '  +6.88%  '

Set-TextValue $ws.Range('D49') '2.033.40'
Set-TextValue $ws.Range('E49') '  +2.10%  '

Set-TextValue $ws.Range('D50') '0.0679'
Set-TextValue $ws.Range('E50') '  +2.84%  '

Set-TextValue $ws.Range('E51') '  +0.10%  '

